$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.186.56"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "2.840.59"
$ws.Range("E3").Value = "  +3.09%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'362.18"
$ws.Range("E5").Value = "  +9.18%  "

$ws.Range("D6").Value = "'115.68"
$ws.Range("E6").Value = "  -1.57%  "

$ws.Range("D7").Value = "'0.551"
$ws.Range("E7").Value = "  +3.36%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.609"
$ws.Range("E9").Value = "  +5.53%  "

$ws.Range("D10").Value = "'42.31"
$ws.Range("E10").Value = "  +1.61%  "

$ws.Range("D11").Value = "'0.0863"
$ws.Range("E11").Value = "  +3.59%  "

$ws.Range("D12").Value = "'20.15"
$ws.Range("E12").Value = "  +0.68%  "

$ws.Range("E13").Value = "  +1.86%  "

$ws.Range("D14").Value = "'7.87"
$ws.Range("E14").Value = "  +3.35%  "

$ws.Range("D15").Value = "3.288.14"
$ws.Range("E15").Value = "  +3.32%  "

$ws.Range("D16").Value = "2.874.85"
$ws.Range("E16").Value = "  +3.54%  "

$ws.Range("D17").Value = "'0.903"
$ws.Range("E17").Value = "  +2.36%  "

$ws.Range("D18").Value = "52.266.89"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("E19").Value = "  +3.22%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'7.32"
$ws.Range("E20").Value = "  +7.11%  "

$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "'3.18"
$ws.Range("E21").Value = "  +5.00%  "

$ws.Range("D22").Value = "0.0₃0994"
$ws.Range("E22").Value = "  +3.30%  "

$ws.Range("D23").Value = "'70.41"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("D24").Value = "'270.92"
$ws.Range("E24").Value = "  -2.88%  "

$ws.Range("E25").Value = "  +6.99%  "

$ws.Range("D26").Value = "'27.16"
$ws.Range("E26").Value = "  +1.25%  "

$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").Value = "'10.31"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").Value = "'0.140"
$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").Value = "'34.56"
$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("D32").Value = "'51.17"
$ws.Range("E32").Value = "  +1.47%  "

$ws.Range("D33").Value = "'5.84"
$ws.Range("E33").Value = "  +4.60%  "

$ws.Range("D34").Value = "'0.0441"
$ws.Range("E34").Value = "  +26.72%  "

$ws.Range("D35").Value = "'0.0836"
$ws.Range("E35").Value = "  +1.63%  "

$ws.Range("D36").Value = "'2.13"
$ws.Range("E36").Value = "  +1.97%  "

$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("D38").Value = "'5.00"
$ws.Range("E38").Value = "  +0.77%  "

$ws.Range("D39").Value = "'3.28"
$ws.Range("E39").Value = "  +2.73%  "

$ws.Range("D40").Value = "'18.64"
$ws.Range("E40").Value = "  -2.22%  "

$ws.Range("E41").Value = "  +7.86%  "

$ws.Range("D42").Value = "'23.74"
$ws.Range("E42").Value = "  +2.46%  "

$ws.Range("E43").Value = "  +2.65%  "

$ws.Range("D44").Value = "'127.88"
$ws.Range("E44").Value = "  -2.38%  "

$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("D46").Value = "'3.40"
$ws.Range("E46").Value = "  +1.92%  "

$ws.Range("D47").Value = "2.074.95"
$ws.Range("E47").Value = "  -1.52%  "

$ws.Range("E48").Value = "  +3.26%  "

$ws.Range("D49").Value = "'0.950"
$ws.Range("E49").Value = "  +8.14%  "

$ws.Range("D50").Value = "'5.61"
$ws.Range("E50").Value = "  +0.69%  "

$ws.Range("E51").Value = "  +0.79%  "
